$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.783.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +8.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.215.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "397.81"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.78"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.07%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.30"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0907"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.27%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.726.56"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.84%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.43%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.09"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.225.98"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.16%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +6.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.71"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "55.687.19"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.36"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.22%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.94%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "302.54"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +13.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.17"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.72%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.19"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.53"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.174"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.64%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.26"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +8.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0493"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "36.12"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.28"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.08"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +22.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.48"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "135.06"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.64%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.92"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.03"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.15"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.286"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.119"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.30"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +46.69%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.49"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.136.63"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0365"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +12.39%  "
